# "Range of updates to report/report data"
#
# Tidies up the categories table on Sheet1:
#  - removes the stray, disconnected formatted cell left behind at C14
#  - fills in the previously-blank "Sub-category" cell for the "Not
#    applicable" row with a single space placeholder
#  - turns on word-wrap for the Main category / Sub-category columns so the
#    long descriptions are readable, and grows the rows whose text now
#    wraps onto two lines
#  - leaves the sheet with columns B:C selected, matching the new
#    authoring view

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the orphaned styled-only cell that used to sit at row 14 (far below
# the real table, which only runs through row 9).
$ws.Rows("14").Delete() | Out-Null

# Row 2 ("Not applicable ...") was missing a Sub-category entry - give it
# the same "blank" placeholder value used elsewhere in the sheet.
$ws.Range("C2").Value = " "

# Wrap text for the Main category (B) and Sub category (C) columns across
# the whole table.
$ws.Range("B1:C9").WrapText = $true

# The descriptions that wrap onto two lines need a taller row to display
# fully.
$ws.Rows("2").RowHeight = 30
$ws.Rows("6").RowHeight = 30
$ws.Rows("9").RowHeight = 30

# Match the updated selection/view state - columns B:C selected.
$ws.Range("B1:C1048576").Select() | Out-Null
